$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph "Line of investigation #1: ..." body
#   "  we researched ..." -> "  We researched ..." (capitalize "we")
#   "corelated" -> "correlated" (fix typo)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("  we researched the relationship", $true, $false, $false, $false, $false, $true, 1, $false, "  We researched the relationship", 2) | Out-Null
$d.Content.Find.Execute("highly corelated", $true, $false, $false, $false, $false, $true, 1, $false, "highly correlated", 2) | Out-Null

# ---------------------------------------------------------------------
# Paragraph "The data was negatively correlated ..."
#   add " quantitative" and a trailing new sentence
# ---------------------------------------------------------------------
$d.Content.Find.Execute("that we found no indication", $true, $false, $false, $false, $false, $true, 1, $false, "that we found no quantitative indication", 2) | Out-Null
$d.Content.Find.Execute("was related to violent crime.", $true, $false, $false, $false, $false, $true, 1, $false, "was related to violent crime. In other words, the number of lead pipes servicing any given population in the United States had no effect on the rates of criminally violent incidents.", 2) | Out-Null

# ---------------------------------------------------------------------
# Merge the "correlation matrix" paragraph with the "However, ... Missouri" paragraph
# ---------------------------------------------------------------------
$pMatrix = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "We relied heavily on the correlation matrix*") {
        $pMatrix = $p
        break
    }
}
$markRange = $d.Range($pMatrix.Range.End - 1, $pMatrix.Range.End)
$markRange.InsertAfter(" ")
$markRange = $d.Range($pMatrix.Range.End - 1, $pMatrix.Range.End)
$markRange.Delete()

# ---------------------------------------------------------------------
# "Line of investigation #2:" paragraph + the 2 following paragraphs get merged
# ---------------------------------------------------------------------
$d.Content.Find.Execute("(EPA data). ", $true, $false, $false, $false, $false, $true, 1, $false, "(EPA data).", 2) | Out-Null
$d.Content.Find.Execute("The notebooks 05-08", $true, $false, $false, $false, $false, $true, 1, $false, " Notebooks 05-08", 2) | Out-Null
$d.Content.Find.Execute("investigate the data.", $true, $false, $false, $false, $false, $true, 1, $false, "investigate the data. ", 2) | Out-Null
$d.Content.Find.Execute("positively correlated but weakly correlated. ", $true, $false, $false, $false, $false, $true, 1, $false, "positively correlated, however, to what would be considered by a “very weak amount”. ", 2) | Out-Null

$pInv2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Line of investigation #2:*") {
        $pInv2 = $p
        break
    }
}
# Merge paragraph with "Notebooks 05-08 ..."
$markRange = $d.Range($pInv2.Range.End - 1, $pInv2.Range.End)
$markRange.Delete()
# Merge paragraph with "The results were positively correlated ..."
$markRange = $d.Range($pInv2.Range.End - 1, $pInv2.Range.End)
$markRange.Delete()
